$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Celecoxib")

# Correct the RxNorm/RxCUI value for Celecoxib: D2 previously held the text
# "C0538927" (a UMLS CUI, which is the wrong kind of identifier for the
# RxNorm row). Replace it with the correct numeric RxCUI, 140587.
$ws.Range("D2").Value = 140587

# Reposition the saved selection/cursor to match where the author left it.
$null = $ws.Range("B11").Select()
